# Fruta / hortaliza, semanal
# Insert a new week of "Feria Lagunitas de Puerto Montt - Manzana" records
# at the top of the data block (rows 979-982), pushing the existing
# historical rows down by 4 (979-1046 -> 983-1050).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows 979:1046 down by 4 rows, preserving them
# unchanged, and opening up 4 fresh rows at 979:982.
$ws.Range("A979:T982").EntireRow.Insert()

# New week date (serial 44585) for the 4 inserted rows.
$newDate = 44585

$ws.Cells.Item(979, 1).Value = 4
$ws.Cells.Item(979, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(979, 3).Value = "Los Lagos"
$ws.Cells.Item(979, 4).Value = $newDate
$ws.Cells.Item(979, 5).Value = 10
$ws.Cells.Item(979, 6).Value = "Fruta"
$ws.Cells.Item(979, 7).Value = 100104
$ws.Cells.Item(979, 8).Value = "Frutos de pepita"
$ws.Cells.Item(979, 9).Value = 100104002
$ws.Cells.Item(979, 10).Value = "Manzana"
$ws.Cells.Item(979, 11).Value = "Fuji royal"
$ws.Cells.Item(979, 12).Value = "Primera"
$ws.Cells.Item(979, 13).Value = 100
$ws.Cells.Item(979, 14).Value = 17000
$ws.Cells.Item(979, 15).Value = 17000
$ws.Cells.Item(979, 16).Value = 17000
$ws.Cells.Item(979, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(979, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(979, 19).Value = 1062
$ws.Cells.Item(979, 20).Value = 16

$ws.Cells.Item(980, 1).Value = 4
$ws.Cells.Item(980, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(980, 3).Value = "Los Lagos"
$ws.Cells.Item(980, 4).Value = $newDate
$ws.Cells.Item(980, 5).Value = 10
$ws.Cells.Item(980, 6).Value = "Fruta"
$ws.Cells.Item(980, 7).Value = 100104
$ws.Cells.Item(980, 8).Value = "Frutos de pepita"
$ws.Cells.Item(980, 9).Value = 100104002
$ws.Cells.Item(980, 10).Value = "Manzana"
$ws.Cells.Item(980, 11).Value = "Fuji royal"
$ws.Cells.Item(980, 12).Value = "Segunda"
$ws.Cells.Item(980, 13).Value = 100
$ws.Cells.Item(980, 14).Value = 14000
$ws.Cells.Item(980, 15).Value = 14000
$ws.Cells.Item(980, 16).Value = 14000
$ws.Cells.Item(980, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(980, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(980, 19).Value = 875
$ws.Cells.Item(980, 20).Value = 16

$ws.Cells.Item(981, 1).Value = 4
$ws.Cells.Item(981, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(981, 3).Value = "Los Lagos"
$ws.Cells.Item(981, 4).Value = $newDate
$ws.Cells.Item(981, 5).Value = 10
$ws.Cells.Item(981, 6).Value = "Fruta"
$ws.Cells.Item(981, 7).Value = 100104
$ws.Cells.Item(981, 8).Value = "Frutos de pepita"
$ws.Cells.Item(981, 9).Value = 100104002
$ws.Cells.Item(981, 10).Value = "Manzana"
$ws.Cells.Item(981, 11).Value = "Granny Smith"
$ws.Cells.Item(981, 12).Value = "Primera"
$ws.Cells.Item(981, 13).Value = 150
$ws.Cells.Item(981, 14).Value = 19000
$ws.Cells.Item(981, 15).Value = 19000
$ws.Cells.Item(981, 16).Value = 19000
$ws.Cells.Item(981, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(981, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(981, 19).Value = 1188
$ws.Cells.Item(981, 20).Value = 16

$ws.Cells.Item(982, 1).Value = 4
$ws.Cells.Item(982, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(982, 3).Value = "Los Lagos"
$ws.Cells.Item(982, 4).Value = $newDate
$ws.Cells.Item(982, 5).Value = 10
$ws.Cells.Item(982, 6).Value = "Fruta"
$ws.Cells.Item(982, 7).Value = 100104
$ws.Cells.Item(982, 8).Value = "Frutos de pepita"
$ws.Cells.Item(982, 9).Value = 100104002
$ws.Cells.Item(982, 10).Value = "Manzana"
$ws.Cells.Item(982, 11).Value = "Granny Smith"
$ws.Cells.Item(982, 12).Value = "Segunda"
$ws.Cells.Item(982, 13).Value = 150
$ws.Cells.Item(982, 14).Value = 15000
$ws.Cells.Item(982, 15).Value = 15000
$ws.Cells.Item(982, 16).Value = 15000
$ws.Cells.Item(982, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(982, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(982, 19).Value = 938
$ws.Cells.Item(982, 20).Value = 16

# Make sure the new date cells carry the same date number format as the
# rest of column D.
$ws.Range("D979:D982").NumberFormat = $ws.Range("D983").NumberFormat
